$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 865.73334
$ws.Range("J28").Value = 1269.8334
$ws.Range("L28").Value = 1269.8334
$ws.Range("N28").Value = -2239.8334
$ws.Range("H74").Value = 12500
$ws.Range("J74").Value = 12000
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -13872
$ws.Range("H76").Value = 7491.5
$ws.Range("I76").Value = 5605.1665
$ws.Range("J76").Value = 8906.25
$ws.Range("K76").Value = 5605.1665
$ws.Range("L76").Value = 8906.25
$ws.Range("M76").Value = -5290.1665
$ws.Range("N76").Value = -9536.25
$ws.Range("H77").Value = 12500
$ws.Range("J77").Value = 12000
$ws.Range("L77").Value = 60000
$ws.Range("N77").Value = -69360
$ws.Range("H79").Value = 7491.5
$ws.Range("I79").Value = 5605.1665
$ws.Range("J79").Value = 8906.25
$ws.Range("K79").Value = 5605.1665
$ws.Range("L79").Value = 8906.25
$ws.Range("M79").Value = -4513.1665
$ws.Range("N79").Value = -11090.25
$ws.Range("H92").Value = 725.913
$ws.Range("I92").Value = 400.13333
$ws.Range("K92").Value = 400.13333
$ws.Range("M92").Value = 847.86667
$ws.Range("H113").Value = 10874.375
$ws.Range("I113").Value = 10997.5
$ws.Range("K113").Value = 10997.5
$ws.Range("M113").Value = -7743.5
$ws.Range("H118").Value = 1529.2222
$ws.Range("I118").Value = 1852.8
$ws.Range("K118").Value = 5558.4
$ws.Range("M118").Value = -3901.4
$ws.Range("H129").Value = 101456.9
$ws.Range("I129").Value = 126137.75
$ws.Range("J129").Value = 2733.5
$ws.Range("K129").Value = 378413.25
$ws.Range("L129").Value = 8200.5
$ws.Range("M129").Value = -373413.25
$ws.Range("N129").Value = -18200.5
$ws.Range("H132").Value = 13575.041
$ws.Range("I132").Value = 2224.2195
$ws.Range("J132").Value = 71748
$ws.Range("K132").Value = 6672.6585
$ws.Range("L132").Value = 215244
$ws.Range("M132").Value = -4142.6585
$ws.Range("N132").Value = -220304
$ws.Range("H137").Value = 4077.3901
$ws.Range("I137").Value = 4421
$ws.Range("K137").Value = 13263
$ws.Range("M137").Value = -10713
$ws.Range("H138").Value = 3002.8313
$ws.Range("I138").Value = 1389.1389
$ws.Range("J138").Value = 4238.851
$ws.Range("K138").Value = 4167.4167
$ws.Range("L138").Value = 12716.553
$ws.Range("M138").Value = 972.5833000000002
$ws.Range("N138").Value = -22996.553
$ws.Range("H141").Value = 5693.7036
$ws.Range("I141").Value = 2942.476
$ws.Range("J141").Value = 15323
$ws.Range("K141").Value = 8827.428
$ws.Range("L141").Value = 45969
$ws.Range("M141").Value = -3647.428
$ws.Range("N141").Value = -56329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17551176
$ws.Range("I32").Value = 19238366
$ws.Range("J32").Value = 4395.8
$ws.Range("K32").Value = 19238366
$ws.Range("L32").Value = 4395.8
$ws.Range("M32").Value = -19238079
$ws.Range("N32").Value = -4969.8
$ws.Range("H45").Value = 6363.375
$ws.Range("I45").Value = 4484.5
$ws.Range("K45").Value = 4484.5
$ws.Range("M45").Value = -4107.5
$ws.Range("H61").Value = 2242.4167
$ws.Range("I61").Value = 2264.4546
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2264.4546
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -2052.4546
$ws.Range("N61").Value = -2424
$ws.Range("H122").Value = 2257.4443
$ws.Range("I122").Value = 1460.8
$ws.Range("K122").Value = 4382.4
$ws.Range("M122").Value = -1932.4
$ws.Range("H132").Value = 1804.4166
$ws.Range("I132").Value = 1755.5758
$ws.Range("J132").Value = 2341.6667
$ws.Range("K132").Value = 5266.7274
$ws.Range("L132").Value = 7025.000100000001
$ws.Range("M132").Value = -2736.7274
$ws.Range("N132").Value = -12085.0001
$ws.Range("H136").Value = 2242.4167
$ws.Range("I136").Value = 2264.4546
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6793.3638
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -4243.3638
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3547.5
$ws.Range("I105").Value = 1913.3334
$ws.Range("K105").Value = 1913.3334
$ws.Range("M105").Value = -166.3334
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1951.5128
$ws.Range("I31").Value = 1847.8182
$ws.Range("K31").Value = 1847.8182
$ws.Range("M31").Value = -1552.8182
$ws.Range("H34").Value = 1951.5128
$ws.Range("I34").Value = 1847.8182
$ws.Range("K34").Value = 1847.8182
$ws.Range("M34").Value = -1645.8182
$ws.Range("H122").Value = 447798.22
$ws.Range("I122").Value = 730525.7
$ws.Range("K122").Value = 2191577.1
$ws.Range("M122").Value = -2189127.1
$ws.Range("H134").Value = 1917.3948
$ws.Range("I134").Value = 1944.2572
$ws.Range("J134").Value = 1604
$ws.Range("K134").Value = 5832.7716
$ws.Range("L134").Value = 4812
$ws.Range("M134").Value = -3297.7716
$ws.Range("N134").Value = -9882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1375.2727
$ws.Range("I5").Value = 1356.5555
$ws.Range("J5").Value = 1459.5
$ws.Range("K5").Value = 4069.6665
$ws.Range("L5").Value = 4378.5
$ws.Range("M5").Value = -3957.6665
$ws.Range("N5").Value = -4602.5
$ws.Range("H95").Value = 6666.6665
$ws.Range("J95").Value = 6666.6665
$ws.Range("L95").Value = 19999.9995
$ws.Range("N95").Value = -24117.9995
$ws.Range("H131").Value = 3424.2812
$ws.Range("I131").Value = 2513
$ws.Range("J131").Value = 3780.8696
$ws.Range("K131").Value = 7539
$ws.Range("L131").Value = 11342.6088
$ws.Range("M131").Value = -2499
$ws.Range("N131").Value = -21422.6088
$ws.Range("H135").Value = 1375.2727
$ws.Range("I135").Value = 1356.5555
$ws.Range("J135").Value = 1459.5
$ws.Range("K135").Value = 12208.9995
$ws.Range("L135").Value = 13135.5
$ws.Range("M135").Value = -9673.9995
$ws.Range("N135").Value = -18205.5
$ws.Range("H137").Value = 2974.1304
$ws.Range("J137").Value = 3270.8333
$ws.Range("L137").Value = 9812.499899999999
$ws.Range("N137").Value = -20012.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 8333484.5
$ws.Range("I2").Value = 49.5
$ws.Range("J2").Value = 16666920
$ws.Range("K2").Value = 49.5
$ws.Range("L2").Value = 16666920
$ws.Range("M2").Value = 63.5
$ws.Range("N2").Value = -16667146
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H80").Value = 114683.63
$ws.Range("I80").Value = 196111
$ws.Range("J80").Value = 16970.8
$ws.Range("K80").Value = 196111
$ws.Range("L80").Value = 16970.8
$ws.Range("M80").Value = -195113
$ws.Range("N80").Value = -18966.8
$ws.Range("H83").Value = 114683.63
$ws.Range("I83").Value = 196111
$ws.Range("J83").Value = 16970.8
$ws.Range("K83").Value = 980555
$ws.Range("L83").Value = 84854
$ws.Range("M83").Value = -975563
$ws.Range("N83").Value = -94838
$ws.Range("H97").Value = 544.1177
$ws.Range("I97").Value = 516.73334
$ws.Range("J97").Value = 749.5
$ws.Range("K97").Value = 516.73334
$ws.Range("L97").Value = 749.5
$ws.Range("M97").Value = -20.73334
$ws.Range("N97").Value = -1741.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 1864.1578
$ws.Range("I132").Value = 1469.9375
$ws.Range("J132").Value = 3966.6667
$ws.Range("K132").Value = 4409.8125
$ws.Range("L132").Value = 11900.0001
$ws.Range("M132").Value = -1879.8125
$ws.Range("N132").Value = -16960.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 6500
$ws.Range("I55").Value = 6000
$ws.Range("J55").Value = 8500
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 8500
$ws.Range("M55").Value = -5827
$ws.Range("N55").Value = -8846
$ws.Range("H93").Value = 5824.5
$ws.Range("I93").Value = 999
$ws.Range("J93").Value = 10650
$ws.Range("K93").Value = 999
$ws.Range("L93").Value = 10650
$ws.Range("M93").Value = 249
$ws.Range("N93").Value = -13146
$ws.Range("H122").Value = 5887.56
$ws.Range("I122").Value = 2714
$ws.Range("J122").Value = 6889.737
$ws.Range("K122").Value = 8142
$ws.Range("L122").Value = 20669.211
$ws.Range("M122").Value = -5692
$ws.Range("N122").Value = -25569.211
$ws.Range("H134").Value = 64000
$ws.Range("J134").Value = 64000
$ws.Range("L134").Value = 64000
$ws.Range("N134").Value = -74140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 27245
$ws.Range("I58").Value = 9995
$ws.Range("K58").Value = 9995
$ws.Range("M58").Value = -9687
$ws.Range("H132").Value = 2089.0344
$ws.Range("I132").Value = 2072.7856
$ws.Range("K132").Value = 6218.3568
$ws.Range("M132").Value = -3688.3568

# Remove cells that must not exist in the final sheet (value becomes structurally absent)
$wb.Worksheets.Item("GSM").Range("N39").ClearContents()
$wb.Worksheets.Item("GSM").Range("M126").ClearContents()
